# Update countries & provincias Spain
# - Refresh "Datos actualizados" timestamp
# - Uzbekistan overtakes Costa de Marfil in the ranking (rows 70/71 swap)
# - Refresh case counters for a handful of countries (Kazajistan, Uzbekistan,
#   Costa de Marfil, Kirguistan, Tailandia, Butan)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 07:29"

# Kazajistan (row 34)
$ws.Range("B34").Value = 56455
$ws.Range("C34").Value = 1708
$ws.Range("D34").Value = 32500
$ws.Range("E34").Value = 23691

# Row 70 now becomes Uzbekistan (was Costa de Marfil) with fresh data
$ws.Range("A70").Value = "Uzbekistan"
$ws.Range("B70").Value = 12206
$ws.Range("C70").Value = 179
$ws.Range("D70").Value = 7530
$ws.Range("E70").Value = 4621
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 55

# Row 71 now becomes Costa de Marfil (was Uzbekistan), keeping its old data
$ws.Range("A71").Value = "Costa de Marfil"
$ws.Range("B71").Value = 12052
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 6080
$ws.Range("E71").Value = 5891
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 81

# Kirguistan (row 73)
$ws.Range("B73").Value = 9910
$ws.Range("C73").Value = 552
$ws.Range("D73").Value = 3236
$ws.Range("E73").Value = 6549

# Tailandia (row 103)
$ws.Range("B103").Value = 3216
$ws.Range("C103").Value = 14
$ws.Range("D103").Value = 3088
$ws.Range("E103").Value = 70

# Butan (row 187)
$ws.Range("B187").Value = 82
$ws.Range("C187").Value = 2
$ws.Range("D187").Value = 57
